$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full updated data block (rows 2-31, columns A-H) in one shot.
# This covers: a newly inserted sensor reading at row 2 (shifting the previously
# recorded rows down by one), plus 10 additional new rows appended at the end.
$data = New-Object 'object[,]' 30,8
$data[0,0] = 0
$data[0,1] = "struggle"
$data[0,2] = -0.01545035839080799
$data[0,3] = 0.2315296996384856
$data[0,4] = -0.5798605158925052
$data[0,5] = -0.0303905457258224
$data[0,6] = -0.02122756652534
$data[0,7] = -0.038026362657547
$data[1,0] = 100
$data[1,1] = "struggle"
$data[1,2] = 0.06989622116088862
$data[1,3] = -0.2621099762618542
$data[1,4] = 0.03679747134447098
$data[1,5] = 0.0122173046693205
$data[1,6] = -0.0598647929728031
$data[1,7] = -0.067195177078247
$data[2,0] = 200
$data[2,1] = "struggle"
$data[2,2] = 0.1602854728698729
$data[2,3] = 0.2304907962679849
$data[2,4] = -0.932947013527152
$data[2,5] = -0.1317941695451736
$data[2,6] = -0.223882108926773
$data[2,7] = -0.0047342055477201
$data[3,0] = 300
$data[3,1] = "struggle"
$data[3,2] = 1.401312828063962
$data[3,3] = -0.6429092288017235
$data[3,4] = -4.608263850212088
$data[3,5] = -0.5534439086914062
$data[3,6] = -2.408488988876343
$data[3,7] = -0.1953241676092147
$data[4,0] = 400
$data[4,1] = "struggle"
$data[4,2] = 2.035469830036161
$data[4,3] = -1.478700906038283
$data[4,4] = -5.709699153900139
$data[4,5] = -0.5908594131469727
$data[4,6] = -4.143346309661865
$data[4,7] = 0.8623890280723572
$data[5,0] = 500
$data[5,1] = "struggle"
$data[5,2] = 0.4589872062206247
$data[5,3] = -0.1242832839488965
$data[5,4] = -1.642685413360593
$data[5,5] = -0.3294090926647186
$data[5,6] = -3.587306261062622
$data[5,7] = 1.745547413825989
$data[6,0] = 600
$data[6,1] = "struggle"
$data[6,2] = -1.821219787001614
$data[6,3] = 1.541860699653622
$data[6,4] = 1.125372886657715
$data[6,5] = -1.008996725082398
$data[6,6] = -2.703842401504517
$data[6,7] = 1.318247199058533
$data[7,0] = 700
$data[7,1] = "struggle"
$data[7,2] = -4.537274837493896
$data[7,3] = -0.381650447845459
$data[7,4] = 0.8837652206420898
$data[7,5] = 0.1534798890352249
$data[7,6] = -0.0916297882795333
$data[7,7] = 0.2032654136419296
$data[8,0] = 800
$data[8,1] = "struggle"
$data[8,2] = -6.140988111495973
$data[8,3] = 2.625270426273349
$data[8,4] = 3.428681612014774
$data[8,5] = 1.313818454742432
$data[8,6] = 3.833332061767578
$data[8,7] = -1.68843150138855
$data[9,0] = 900
$data[9,1] = "struggle"
$data[9,2] = -6.243353843688965
$data[9,3] = 0.4006794318556786
$data[9,4] = -2.078256130218506
$data[9,5] = 0.2997821271419525
$data[9,6] = 5.671730995178223
$data[9,7] = -1.668731093406677
$data[10,0] = 1000
$data[10,1] = "struggle"
$data[10,2] = 1.525731027126295
$data[10,3] = -1.101333875209092
$data[10,4] = -3.408103108406068
$data[10,5] = -0.6543893814086914
$data[10,6] = 4.36127233505249
$data[10,7] = -0.7669413089752197
$data[11,0] = 1100
$data[11,1] = "struggle"
$data[11,2] = 2.08496618270874
$data[11,3] = -1.613636314868927
$data[11,4] = -2.943879842758179
$data[11,5] = -0.5648976564407349
$data[11,6] = 2.72186279296875
$data[11,7] = 1.022893905639648
$data[12,0] = 1200
$data[12,1] = "struggle"
$data[12,2] = -2.373238265514372
$data[12,3] = -0.1506620794534765
$data[12,4] = -2.15161240100861
$data[12,5] = -0.7005097270011902
$data[12,6] = 1.526399493217468
$data[12,7] = 0.2272418737411499
$data[13,0] = 1300
$data[13,1] = "struggle"
$data[13,2] = -3.501610040664678
$data[13,3] = 2.79044055938721
$data[13,4] = -2.906912446022055
$data[13,5] = 0.0536034256219863
$data[13,6] = -1.495856285095215
$data[13,7] = -0.2338086664676666
$data[14,0] = 1400
$data[14,1] = "struggle"
$data[14,2] = -4.058286607265482
$data[14,3] = 2.954759478569038
$data[14,4] = -7.582780838012708
$data[14,5] = 0.1302670091390609
$data[14,6] = -5.946773052215576
$data[14,7] = -1.086271166801453
$data[15,0] = 1500
$data[15,1] = "struggle"
$data[15,2] = 2.394027709960937
$data[15,3] = -1.486701488494873
$data[15,4] = 0.1956815719604492
$data[15,5] = -0.3718642294406891
$data[15,6] = -4.293160915374756
$data[15,7] = 1.783268332481384
$data[16,0] = 1600
$data[16,1] = "struggle"
$data[16,2] = -2.034749180078514
$data[16,3] = -1.211341693997383
$data[16,4] = -0.07141649723053026
$data[16,5] = 0.1963931769132614
$data[16,6] = -4.29804801940918
$data[16,7] = 1.745852828025818
$data[17,0] = 1700
$data[17,1] = "struggle"
$data[17,2] = -4.999244511127475
$data[17,3] = -0.394761115312575
$data[17,4] = -0.1581084728240968
$data[17,5] = 0.204487144947052
$data[17,6] = -2.893821477890014
$data[17,7] = 1.158658623695374
$data[18,0] = 1800
$data[18,1] = "struggle"
$data[18,2] = -6.522920310497289
$data[18,3] = 0.375426143407825
$data[18,4] = -0.1347707509994504
$data[18,5] = 0.1171334087848663
$data[18,6] = 0.94042706489563
$data[18,7] = -0.3061962127685547
$data[19,0] = 1900
$data[19,1] = "struggle"
$data[19,2] = -6.107214450836182
$data[19,3] = 0.1423146724700928
$data[19,4] = 0.9927992820739744
$data[19,5] = 0.337044894695282
$data[19,6] = 3.794236660003662
$data[19,7] = -1.184467673301697
$data[20,0] = 2000
$data[20,1] = "struggle"
$data[20,2] = -4.145634770393367
$data[20,3] = -0.62879066169262
$data[20,4] = 1.90212270617485
$data[20,5] = -0.3572034537792206
$data[20,6] = 5.656154155731201
$data[20,7] = -1.049466490745544
$data[21,0] = 2100
$data[21,1] = "struggle"
$data[21,2] = -1.618811368942268
$data[21,3] = -0.6589505374431646
$data[21,4] = 0.2403407692909387
$data[21,5] = -0.0070249503478407
$data[21,6] = 4.270253658294678
$data[21,7] = -0.0296269636601209
$data[22,0] = 2200
$data[22,1] = "struggle"
$data[22,2] = -0.08240008354186718
$data[22,3] = 0.4783504903316493
$data[22,4] = -3.809414207935333
$data[22,5] = 0.4137084782123565
$data[22,6] = 2.936276435852051
$data[22,7] = 0.2823724448680877
$data[23,0] = 2300
$data[23,1] = "struggle"
$data[23,2] = -3.95973014831543
$data[23,3] = 0.9762580394744872
$data[23,4] = -4.069071769714356
$data[23,5] = 0.0429132841527462
$data[23,6] = 1.122159481048584
$data[23,7] = 0.1867720484733581
$data[24,0] = 2400
$data[24,1] = "struggle"
$data[24,2] = -3.901577949523926
$data[24,3] = 1.771272063255311
$data[24,4] = -2.484678864479062
$data[24,5] = 0.0647517144680023
$data[24,6] = -1.842216849327088
$data[24,7] = -0.6108652353286743
$data[25,0] = 2500
$data[25,1] = "struggle"
$data[25,2] = -1.957046031951897
$data[25,3] = -0.6577051877975557
$data[25,4] = -7.9572014808655
$data[25,5] = 0.0862847194075584
$data[25,6] = -5.713422775268555
$data[25,7] = -1.346194267272949
$data[26,0] = 2600
$data[26,1] = "struggle"
$data[26,2] = 0.0412573218345611
$data[26,3] = -3.739429324865336
$data[26,4] = -3.584903955459609
$data[26,5] = -0.1818851232528686
$data[26,6] = -4.851491928100586
$data[26,7] = 1.392772793769836
$data[27,0] = 2700
$data[27,1] = "struggle"
$data[27,2] = -2.076164960861222
$data[27,3] = -2.966795355081547
$data[27,4] = -0.5832877159118695
$data[27,5] = -0.3181080818176269
$data[27,6] = -3.869678497314453
$data[27,7] = 0.9886853694915771
$data[28,0] = 2800
$data[28,1] = "struggle"
$data[28,2] = -5.426012933254244
$data[28,3] = -0.3007338047027568
$data[28,4] = 0.185311913490301
$data[28,5] = 0.1050688251852989
$data[28,6] = -2.216677188873291
$data[28,7] = 0.3729332387447357
$data[29,0] = 2900
$data[29,1] = "struggle"
$data[29,2] = -6.20224690437317
$data[29,3] = 0.8901370018720657
$data[29,4] = 1.72858691215514
$data[29,5] = 0.1996002197265625
$data[29,6] = 1.434922456741333
$data[29,7] = -0.2237294018268585
$ws.Range("A2:H31").Value2 = $data

Write-Host "Dimension now:" $ws.UsedRange.Address()
